# "update on Feb 14" - add newly confirmed cases (rows 60-68, case 59-67),
# two trailing blank rows (69-70), and fill in previously-missing
# Status/Visited/SymtomDate/DischargeDate/Stay details for several
# existing cases (rows 31, 46, 55-59) now that more info is available.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- existing rows: fill previously blank cells now that data is known ----

# Case 30 (row 31) - discharged, symptom date now known
$ws.Range("K31").Value = "Discharged"
$ws.Range("N31").Value = "Feb-14"

# Case 45 (row 46) - discharged, symptom date now known
$ws.Range("K46").Value = "Discharged"
$ws.Range("N46").Value = "Feb-14"

# Case 54 (row 55) - corrected coordinates/age, stay/visited/symptom date
$ws.Range("B55").Value = 1.375907
$ws.Range("C55").Value = 103.854193
$ws.Range("E55").Value = 55
$ws.Range("H55").Value = "Ang Mo Kio Avenue 5"
$ws.Range("I55").Value = "Grace Assembly of God (Tanglin)"
$ws.Range("M55").Value = "Feb-10"

# Case 55 (row 56) - visited/symptom date
$ws.Range("I56").Value = "Pulau Bukom, Church of Christ the King (2221 Ang Mo Kio Avenue 8), GP Clinic"
$ws.Range("M56").Value = "Jan-30"

# Case 56 (row 57) - visited/linked-to/symptom date
$ws.Range("I57").Value = "Seletar Aerospace Heights construction site, Tan Tock Seng Hospital"
$ws.Range("J57").Value = "47"
$ws.Range("M57").Value = "Feb-12"

# Case 57 (row 58) - corrected coordinates, stay/symptom date
$ws.Range("B58").Value = 1.3868290000000001
$ws.Range("C58").Value = 103.761758
$ws.Range("H58").Value = "Senja Road"
$ws.Range("M58").Value = "Feb-11"

# Case 58 (row 59) - corrected coordinates, stay/symptom date
$ws.Range("B59").Value = 1.3467119999999999
$ws.Range("C59").Value = 103.88141899999999
$ws.Range("H59").Value = "Jalan Kelichap"
$ws.Range("M59").Value = "Feb-10"

# ---- new rows 60-70, copy formatting from row 59 (last existing data row) ----

$ws.Range("A59:O59").Copy()
$ws.Range("A60:O68").PasteSpecial(-4122)
$ws.Range("A59:O59").Copy()
$ws.Range("A69:O70").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Case 59 (row 60)
$ws.Range("A60").Value = 59
$ws.Range("B60").Value = 1.302775
$ws.Range("C60").Value = 103.891108
$ws.Range("D60").Value = "Feb-13"
$ws.Range("E60").Value = 61
$ws.Range("F60").Value = "Male"
$ws.Range("G60").Value = "Singapore"
$ws.Range("H60").Value = "Wilkinson Road"
$ws.Range("I60").Value = "Farrer Park Hospital’s emergency care clinic, Mount Elizabeth Hospital emergency department"
$ws.Range("J60").Value = ""
$ws.Range("K60").Value = ""
$ws.Range("L60").Value = "Singapore"
$ws.Range("M60").Value = "Feb-07"
$ws.Range("N60").Value = ""
$ws.Range("O60").Value = ""

# Case 60 (row 61)
$ws.Range("A61").Value = 60
$ws.Range("B61").Value = 1.369327
$ws.Range("C61").Value = 103.85627599999999
$ws.Range("D61").Value = "Feb-13"
$ws.Range("E61").Value = 51
$ws.Range("F61").Value = "Female"
$ws.Range("G61").Value = "Singapore"
$ws.Range("H61").Value = "Ang Mo Kio Avenue 3"
$ws.Range("I61").Value = "Grace Assembly of God (Tanglin), GP Clinic"
$ws.Range("J61").Value = ""
$ws.Range("K61").Value = ""
$ws.Range("L61").Value = "Singapore"
$ws.Range("M61").Value = "Feb-08"
$ws.Range("N61").Value = ""
$ws.Range("O61").Value = "Grace Assembly of God"

# Case 61 (row 62)
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = 1.341642
$ws.Range("C62").Value = 103.70378700000001
$ws.Range("D62").Value = "Feb-13"
$ws.Range("E62").Value = 57
$ws.Range("F62").Value = "Male"
$ws.Range("G62").Value = "Singapore"
$ws.Range("H62").Value = "Jurong West Street 64"
$ws.Range("I62").Value = "Grace Assembly of God, Legacy Office Supplies Pte Ltd (56 Senang Crescent), GP clinic, Ng Teng Fong General Hospital, Pioneer Polyclinic, National University Hospital"
$ws.Range("J62").Value = ""
$ws.Range("K62").Value = ""
$ws.Range("L62").Value = "Singapore"
$ws.Range("M62").Value = "Feb-06"
$ws.Range("N62").Value = ""
$ws.Range("O62").Value = "Grace Assembly of God"

# Case 62 (row 63)
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = 1.2763329999999999
$ws.Range("C63").Value = 103.843384
$ws.Range("D63").Value = "Feb-13"
$ws.Range("E63").Value = 44
$ws.Range("F63").Value = "Female"
$ws.Range("G63").Value = "Singapore"
$ws.Range("H63").Value = "Tanjong Pagar Road"
$ws.Range("I63").Value = "Grace Assembly of God, GP clinic, Singapore General Hospital "
$ws.Range("J63").Value = ""
$ws.Range("K63").Value = ""
$ws.Range("L63").Value = "Singapore"
$ws.Range("M63").Value = "Feb-09"
$ws.Range("N63").Value = ""
$ws.Range("O63").Value = "Grace Assembly of God"

# Case 63 (row 64)
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = 1.2942899999999999
$ws.Range("C64").Value = 103.81654899999999
$ws.Range("D64").Value = "Feb-14"
$ws.Range("E64").Value = 54
$ws.Range("F64").Value = "Female"
$ws.Range("G64").Value = "Singapore"
$ws.Range("H64").Value = ""
$ws.Range("I64").Value = "Grace Assembly of God"
$ws.Range("J64").Value = ""
$ws.Range("K64").Value = ""
$ws.Range("L64").Value = "Singapore"
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = ""
$ws.Range("O64").Value = "Grace Assembly of God"

# Case 64 (row 65)
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = 1.3220970000000001
$ws.Range("C65").Value = 103.847273
$ws.Range("D65").Value = "Feb-14"
$ws.Range("E65").Value = 50
$ws.Range("F65").Value = "Male"
$ws.Range("G65").Value = "Singapore"
$ws.Range("H65").Value = ""
$ws.Range("I65").Value = ""
$ws.Range("J65").Value = ""
$ws.Range("K65").Value = ""
$ws.Range("L65").Value = "Singapore"
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = ""
$ws.Range("O65").Value = ""

# Case 65 (row 66)
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = 1.353793
$ws.Range("C66").Value = 103.86065000000001
$ws.Range("D66").Value = "Feb-14"
$ws.Range("E66").Value = 61
$ws.Range("F66").Value = "Female"
$ws.Range("G66").Value = "Singapore"
$ws.Range("H66").Value = "Mei Hwan Drive"
$ws.Range("I66").Value = ""
$ws.Range("J66").Value = "55"
$ws.Range("K66").Value = ""
$ws.Range("L66").Value = "Singapore"
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = ""
$ws.Range("O66").Value = ""

# Case 66 (row 67)
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 1.294427
$ws.Range("C67").Value = 103.816418
$ws.Range("D67").Value = "Feb-14"
$ws.Range("E67").Value = 28
$ws.Range("F67").Value = "Male"
$ws.Range("G67").Value = "Singapore"
$ws.Range("H67").Value = ""
$ws.Range("I67").Value = "Grace Assembly of God"
$ws.Range("J67").Value = ""
$ws.Range("K67").Value = ""
$ws.Range("L67").Value = "Singapore"
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = ""
$ws.Range("O67").Value = "Grace Assembly of God"

# Case 67 (row 68)
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = 1.3405910000000001
$ws.Range("C68").Value = 103.703937
$ws.Range("D68").Value = "Feb-14"
$ws.Range("E68").Value = 56
$ws.Range("F68").Value = "Female"
$ws.Range("G68").Value = "Singapore"
$ws.Range("H68").Value = "Jurong West Street 64"
$ws.Range("I68").Value = "Grace Assembly of God"
$ws.Range("J68").Value = "61"
$ws.Range("K68").Value = ""
$ws.Range("L68").Value = "Singapore"
$ws.Range("M68").Value = ""
$ws.Range("N68").Value = ""
$ws.Range("O68").Value = "Grace Assembly of God"

# Rows 69-70: blank spacer rows, only column A keeps the row's style (no values)
$ws.Range("B69:O69").ClearContents()
$ws.Range("B70:O70").ClearContents()

$ws.Range("I70").Select()
